$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.410.41'
$ws.Range("E2").Value = '  +2.54%  '

$ws.Range("D3").Value = '2.109.73'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '''345.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("D7").Value = '''0.5237'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.31%  '

$ws.Range("D8").Value = '''0.4449'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.20%  '

$ws.Range("D9").Value = '''54.78'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.62%  '

$ws.Range("D10").Value = '''0.09395'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.93%  '

$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("D12").Value = '''24.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '

$ws.Range("D13").Value = '''8.708'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.38%  '

$ws.Range("D14").Value = '''6.950'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '2.056.33'
$ws.Range("E15").Value = '  -2.31%  '

$ws.Range("D16").Value = '''101.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.43%  '

$ws.Range("D17").Value = '''0.00001164'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").Value = '''1.007'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("D20").Value = '''0.06725'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("D21").Value = '''6.340'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.80%  '

$ws.Range("D22").Value = '''1.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").Value = '30.443.76'
$ws.Range("E23").Value = '  +2.46%  '

$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("D25").Value = '''2.301'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.38%  '

$ws.Range("D26").Value = '''22.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '

$ws.Range("D27").Value = '''162.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").Value = '''2.532'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.65%  '

$ws.Range("D29").Value = '''134.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").Value = '''1.154'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.36%  '

$ws.Range("D31").Value = '''1.740'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.62%  '

$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("D33").Value = '''6.828'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.19%  '

$ws.Range("D34").Value = '''6.271'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.88%  '

$ws.Range("D35").Value = '''3.919'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '

$ws.Range("D36").Value = '''10.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.13%  '

$ws.Range("D37").Value = '''0.02627'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.25%  '

$ws.Range("D38").Value = '''0.06793'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.91%  '

$ws.Range("D39").Value = '''0.7055'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.06%  '

$ws.Range("D40").Value = '''1.363'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.33%  '

$ws.Range("D42").Value = '''0.2225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").Value = '''0.6862'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.73%  '

$ws.Range("D44").Value = '''14.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.19%  '

$ws.Range("D45").Value = '''2.363'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.23%  '

$ws.Range("D46").Value = '''1.006'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = '''1.380'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +18.77%  '

$ws.Range("D48").Value = '''3.646'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("D49").Value = '''0.00000000349'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.79%  '

$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''1.207'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.20%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '''1.220'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
